$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F "想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1691
$ws1.Range("F6").Value = 3276
$ws1.Range("F7").Value = 887
$ws1.Range("F8").Value = 2094
$ws1.Range("F9").Value = 2007
$ws1.Range("F10").Value = 1038
$ws1.Range("F18").Value = 97
$ws1.Range("F19").Value = 1468
$ws1.Range("F20").Value = 547
$ws1.Range("F23").Value = 11827
$ws1.Range("F24").Value = 11850
$ws1.Range("F25").Value = 867
$ws1.Range("F29").Value = 471

# Sheet "全部类型" updates (same events repeated, column F "想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1691
$ws4.Range("F8").Value = 3276
$ws4.Range("F9").Value = 887
$ws4.Range("F10").Value = 2094
$ws4.Range("F11").Value = 2007
$ws4.Range("F12").Value = 1038
$ws4.Range("F22").Value = 97
$ws4.Range("F23").Value = 1468
$ws4.Range("F24").Value = 547
$ws4.Range("F27").Value = 11827
$ws4.Range("F28").Value = 11850
$ws4.Range("F29").Value = 867
$ws4.Range("F35").Value = 471
